# Refactor to OOP: drop the GICS_Compliance sheet and zero-out the stale
# "20.16" dummy Total Assets figure (and everything computed from it)
# across the remaining compliance worksheets.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1. Remove the GICS_Compliance worksheet entirely ----------------------
$gics = $wb.Worksheets.Item("GICS_Compliance")
$gics.Delete()

# --- 2. 40Act_Diversification -----------------------------------------------
$ws = $wb.Worksheets.Item("40Act_Diversification")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = "(JNJ, 0, 0.00%, 0.00%), (JNJ, 0, 0.00%, 0.00%), (MRK, 0, 0.00%, 0.00%), (MRK, 0, 0.00%, 0.00%), (CSCO, 0, 0.00%, 0.00%), (CSCO, 0, 0.00%, 0.00%), (KO, 0, 0.00%, 0.00%), (MCD, 0, 0.00%, 0.00%), (MCD, 0, 0.00%, 0.00%), (CVX, 0, 0.00%, 0.00%), (KO, 0, 0.00%, 0.00%), (CVX, 0, 0.00%, 0.00%), (PG, 0, 0.00%, 0.00%), (PG, 0, 0.00%, 0.00%), (AMGN, 0, 0.00%, 0.00%), (AMGN, 0, 0.00%, 0.00%), (VZ, 0, 0.00%, 0.00%), (VZ, 0, 0.00%, 0.00%), (JNJ, 0, 0.00%, 0.00%), (CSCO, 0, 0.00%, 0.00%), (KO, 0, 0.00%, 0.00%), (MCD, 0, 0.00%, 0.00%), (CVX, 0, 0.00%, 0.00%), (AMGN, 0, 0.00%, 0.00%), (VZ, 0, 0.00%, 0.00%), (PG, 0, 0.00%, 0.00%), (MRK, 0, 0.00%, 0.00%)"
$ws.Range("T2").Value = 0
# Note: COM ColumnWidth is offset by ~0.83 from the raw OOXML <col width>
# units on this sheet's default font, so 12 there needs 11.17 here.
$ws.Columns.Item(20).ColumnWidth = 11.17

# --- 3. IRS_Diversification ---------------------------------------------
$ws = $wb.Worksheets.Item("IRS_Diversification")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 55442243.79
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = "(AMGN, 0.00%), (AMGN, 0.00%), (VZ, 0.00%), (VZ, 0.00%), (PG, 0.00%), (PG, 0.00%), (PG, 0.00%), (MRK, 0.00%), (MRK, 0.00%), (MRK, 0.00%), (MCD, 0.00%), (MCD, 0.00%), (MCD, 0.00%), (KO, 0.00%), (KO, 0.00%), (KO, 0.00%), (JNJ, 0.00%), (JNJ, 0.00%), (JNJ, 0.00%), (IBM, 0.00%), (IBM, 0.00%), (IBM, 0.00%), (CVX, 0.00%), (CVX, 0.00%), (CVX, 0.00%), (CSCO, 0.00%), (CSCO, 0.00%), (CSCO, 0.00%), (AMGN, 0.00%), (VZ, 0.00%)"
$ws.Range("O2").Value = "N/A (0.00%)"
$ws.Range("P2").Value = "N/A (0.00%)"

# --- 4. Illiquid -------------------------------------------------------
$ws = $wb.Worksheets.Item("Illiquid")
$ws.Range("C2").Value = 0
$ws.Range("F2").Value = 0

# --- 5. 12d1_Other_Investment_Companies ---------------------------------
$ws = $wb.Worksheets.Item("12d1_Other_Investment_Companies")
$ws.Range("C2").Value = 0

# --- 6. 12d2_Insurance_Companies ----------------------------------------
$ws = $wb.Worksheets.Item("12d2_Insurance_Companies")
$ws.Range("D2").Value = 0

# --- 7. 12d3_Securities_Business ----------------------------------------
$ws = $wb.Worksheets.Item("12d3_Securities_Business")
$ws.Range("J2").Value = 0
